$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3393
$ws.Range("I43").Value = 2251.75
$ws.Range("K43").Value = 2251.75
$ws.Range("M43").Value = -2182.75
$ws.Range("H51").Value = 3740
$ws.Range("J51").Value = 2752.1667
$ws.Range("L51").Value = 2752.1667
$ws.Range("N51").Value = -3720.1667
$ws.Range("H98").Value = 2051.8462
$ws.Range("J98").Value = 1856.2
$ws.Range("L98").Value = 1856.2
$ws.Range("N98").Value = -4852.2
$ws.Range("H116").Value = 9485.75
$ws.Range("I116").Value = 12472.75
$ws.Range("K116").Value = 12472.75
$ws.Range("M116").Value = -9030.75
$ws.Range("H122").Value = 2051.8462
$ws.Range("J122").Value = 1856.2
$ws.Range("L122").Value = 5568.6
$ws.Range("N122").Value = -10468.6
$ws.Range("H132").Value = 1282.5161
$ws.Range("I132").Value = 1158.7667
$ws.Range("K132").Value = 3476.300099999999
$ws.Range("M132").Value = -946.3000999999995
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("N140").ClearContents()

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H43").Value = 11900
$ws.Range("J43").Value = 11900
$ws.Range("L43").Value = 11900
$ws.Range("N43").Value = -12526
$ws.Range("H132").Value = 3756.4
$ws.Range("I132").Value = 3719.8235
$ws.Range("K132").Value = 11159.4705
$ws.Range("M132").Value = -8629.470499999999

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 5994.3335
$ws.Range("J86").Value = 10000
$ws.Range("L86").Value = 10000
$ws.Range("N86").Value = -12246
$ws.Range("H89").Value = 5994.3335
$ws.Range("J89").Value = 10000
$ws.Range("L89").Value = 50000
$ws.Range("N89").Value = -61232
$ws.Range("H107").Value = 2006.4857
$ws.Range("I107").Value = 2073.8215
$ws.Range("J107").Value = 1737.1428
$ws.Range("K107").Value = 2073.8215
$ws.Range("L107").Value = 1737.1428
$ws.Range("M107").Value = -153.8215
$ws.Range("N107").Value = -5577.1428

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1909.5714
$ws.Range("I16").Value = 1961.3334
$ws.Range("J16").Value = 1599
$ws.Range("K16").Value = 1961.3334
$ws.Range("L16").Value = 1599
$ws.Range("M16").Value = -1674.3334
$ws.Range("N16").Value = -2173
$ws.Range("H82").Value = 30000
$ws.Range("J82").Value = 30000
$ws.Range("L82").Value = 30000
$ws.Range("N82").Value = -30722
$ws.Range("H85").Value = 30000
$ws.Range("J85").Value = 30000
$ws.Range("L85").Value = 30000
$ws.Range("N85").Value = -32496
$ws.Range("H107").Value = 2120.5833
$ws.Range("I107").Value = 911
$ws.Range("J107").Value = 2362.5
$ws.Range("K107").Value = 911
$ws.Range("L107").Value = 2362.5
$ws.Range("M107").Value = 1009
$ws.Range("N107").Value = -6202.5
$ws.Range("H113").Value = 1909.5714
$ws.Range("I113").Value = 1961.3334
$ws.Range("J113").Value = 1599
$ws.Range("K113").Value = 1961.3334
$ws.Range("L113").Value = 1599
$ws.Range("M113").Value = 208.6666
$ws.Range("N113").Value = -5939
$ws.Range("H132").Value = 3902.75
$ws.Range("I132").Value = 3902.75
$ws.Range("K132").Value = 11708.25
$ws.Range("M132").Value = -9178.25
$ws.Range("H133").Value = 29888
$ws.Range("J133").Value = 29888
$ws.Range("L133").Value = 29888
$ws.Range("N133").Value = -34948

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 996
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H64").Value = 400
$ws.Range("I64").Value = 400
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 1200
$ws.Range("M64").Value = -930
$ws.Range("N64").Value = 0
$ws.Range("H67").Value = 400
$ws.Range("I67").Value = 400
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 1200
$ws.Range("M67").Value = -264
$ws.Range("N67").Value = 0
$ws.Range("H116").Value = 433.33334
$ws.Range("I116").Value = 400
$ws.Range("K116").Value = 1200
$ws.Range("M116").Value = 2242
$ws.Range("H131").Value = 1719.1578
$ws.Range("I131").Value = 1188.375
$ws.Range("J131").Value = 1860.7
$ws.Range("K131").Value = 3565.125
$ws.Range("L131").Value = 5582.1
$ws.Range("M131").Value = 1474.875
$ws.Range("N131").Value = -15662.1
$ws.Range("H133").Value = 4297.923
$ws.Range("I133").Value = 4239.4165
$ws.Range("K133").Value = 12718.2495
$ws.Range("M133").Value = -7658.249500000002
$ws.Range("H135").Value = 996
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H140").Value = 2180.2856
$ws.Range("I140").Value = 2180.2856
$ws.Range("K140").Value = 6540.8568
$ws.Range("M140").Value = -1360.8568
$ws.Range("H141").Value = 4221.778
$ws.Range("J141").Value = 4812.5
$ws.Range("L141").Value = 14437.5
$ws.Range("N141").Value = -24797.5

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1614.3478
$ws.Range("I102").Value = 1459.579
$ws.Range("K102").Value = 1459.579
$ws.Range("M102").Value = 162.421
$ws.Range("H103").Value = 96383.5
$ws.Range("J103").Value = 96383.5
$ws.Range("L103").Value = 96383.5
$ws.Range("N103").Value = -98727.5
$ws.Range("H107").Value = 1217.5
$ws.Range("J107").Value = 1440
$ws.Range("L107").Value = 1440
$ws.Range("N107").Value = -5280
$ws.Range("H111").Value = 96763
$ws.Range("J111").Value = 96763
$ws.Range("L111").Value = 96763
$ws.Range("N111").Value = -102897
$ws.Range("H113").Value = 87269.86
$ws.Range("I113").Value = 55455
$ws.Range("J113").Value = 99995.8
$ws.Range("K113").Value = 55455
$ws.Range("L113").Value = 99995.8
$ws.Range("M113").Value = -53285
$ws.Range("N113").Value = -104335.8
$ws.Range("H122").Value = 4177.6
$ws.Range("I122").Value = 4177.6
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12532.8
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 2685.762
$ws.Range("J126").Value = 2876.7334
$ws.Range("L126").Value = 8630.200199999999
$ws.Range("N126").Value = -13570.2002
$ws.Range("H132").Value = 4012.1667
$ws.Range("I132").Value = 4012.1667
$ws.Range("K132").Value = 12036.5001
$ws.Range("M132").Value = -9506.500100000001

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H46").Value = 3650.4783
$ws.Range("J46").Value = 3650.4783
$ws.Range("L46").Value = 3650.4783
$ws.Range("N46").Value = -4026.4783
$ws.Range("H55").Value = 1616.8148
$ws.Range("I55").Value = 1601.7646
$ws.Range("K55").Value = 1601.7646
$ws.Range("M55").Value = -1428.7646
$ws.Range("H61").Value = 3500
$ws.Range("I61").Value = 3500
$ws.Range("J61").Value = 3500
$ws.Range("K61").Value = 3500
$ws.Range("L61").Value = 3500
$ws.Range("M61").Value = -3298
$ws.Range("N61").Value = -3904
$ws.Range("H68").Value = 1060.4546
$ws.Range("I68").Value = 682.8
$ws.Range("K68").Value = 682.8
$ws.Range("M68").Value = 66.20000000000005
$ws.Range("H71").Value = 1060.4546
$ws.Range("I71").Value = 682.8
$ws.Range("K71").Value = 3414
$ws.Range("M71").Value = 330
$ws.Range("H113").Value = 3500
$ws.Range("I113").Value = 3500
$ws.Range("J113").Value = 3500
$ws.Range("K113").Value = 3500
$ws.Range("L113").Value = 3500
$ws.Range("M113").Value = -1330
$ws.Range("N113").Value = -7840
$ws.Range("H119").Value = 99929
$ws.Range("J119").Value = 99929
$ws.Range("L119").Value = 99929
$ws.Range("N119").Value = -109605
$ws.Range("H122").Value = 9969.467000000001
$ws.Range("I122").Value = 10328.81
$ws.Range("J122").Value = 9131
$ws.Range("K122").Value = 30986.43
$ws.Range("L122").Value = 27393
$ws.Range("M122").Value = -28536.43
$ws.Range("N122").Value = -32293

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 391.16666
$ws.Range("I113").Value = 361.75
$ws.Range("K113").Value = 1085.25
$ws.Range("M113").Value = 1084.75
$ws.Range("H132").Value = 2543.7297
$ws.Range("I132").Value = 2096.0303
$ws.Range("K132").Value = 6288.090899999999
$ws.Range("M132").Value = -3758.090899999999
